$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7195.909
$ws.Range("J69").Value = 7195.909
$ws.Range("L69").Value = 21587.727
$ws.Range("N69").Value = -23335.727
$ws.Range("H72").Value = 7195.909
$ws.Range("J72").Value = 7195.909
$ws.Range("L72").Value = 64763.181
$ws.Range("N72").Value = -73499.181
$ws.Range("H74").Value = 9914.235000000001
$ws.Range("I74").Value = 9914.235000000001
$ws.Range("K74").Value = 9914.235000000001
$ws.Range("M74").Value = -8978.235000000001
$ws.Range("H77").Value = 9914.235000000001
$ws.Range("I77").Value = 9914.235000000001
$ws.Range("K77").Value = 49571.175
$ws.Range("M77").Value = -44891.175
$ws.Range("H94").Value = 675.3333
$ws.Range("I94").Value = 675.3333
$ws.Range("K94").Value = 675.3333
$ws.Range("M94").Value = -224.3333
$ws.Range("H112").Value = 2567.3547
$ws.Range("J112").Value = 2677.1365
$ws.Range("L112").Value = 8031.4095
$ws.Range("N112").Value = -10247.4095
$ws.Range("H133").Value = 72958.5
$ws.Range("J133").Value = 72958.5
$ws.Range("L133").Value = 72958.5
$ws.Range("N133").Value = -83078.5
$ws.Range("H136").Value = 99703.836
$ws.Range("J136").Value = 99703.836
$ws.Range("L136").Value = 99703.836
$ws.Range("N136").Value = -109903.836
$ws.Range("H137").Value = 10898.789
$ws.Range("I137").Value = 2371
$ws.Range("J137").Value = 22624.5
$ws.Range("K137").Value = 7113
$ws.Range("L137").Value = 67873.5
$ws.Range("M137").Value = -4563
$ws.Range("N137").Value = -72973.5
$ws.Range("H138").Value = 2324.3
$ws.Range("J138").Value = 2586.5715
$ws.Range("L138").Value = 7759.7145
$ws.Range("N138").Value = -18039.7145
$ws.Range("H139").Value = 64299.2
$ws.Range("J139").Value = 64299.2
$ws.Range("L139").Value = 64299.2
$ws.Range("N139").Value = -74579.2
$ws.Range("H141").Value = 2291.4211
$ws.Range("I141").Value = 2385.6667
$ws.Range("K141").Value = 7157.000100000001
$ws.Range("M141").Value = -1977.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 151467.94
$ws.Range("I32").Value = 163988.77
$ws.Range("J32").Value = 19999.334
$ws.Range("K32").Value = 163988.77
$ws.Range("L32").Value = 19999.334
$ws.Range("M32").Value = -163701.77
$ws.Range("N32").Value = -20573.334
$ws.Range("H63").Value = 3473.2
$ws.Range("I63").Value = 3108.1667
$ws.Range("K63").Value = 3108.1667
$ws.Range("M63").Value = -2422.1667
$ws.Range("H66").Value = 3473.2
$ws.Range("I66").Value = 3108.1667
$ws.Range("K66").Value = 15540.8335
$ws.Range("M66").Value = -12108.8335
$ws.Range("H101").Value = 740000
$ws.Range("J101").Value = 740000
$ws.Range("L101").Value = 740000
$ws.Range("N101").Value = -746490
$ws.Range("H102").Value = 9747.733
$ws.Range("I102").Value = 9747.733
$ws.Range("K102").Value = 9747.733
$ws.Range("M102").Value = -8125.733
$ws.Range("H122").Value = 26317146
$ws.Range("I122").Value = 35715320
$ws.Range("J122").Value = 2260.6
$ws.Range("K122").Value = 107145960
$ws.Range("L122").Value = 6781.799999999999
$ws.Range("M122").Value = -107143510
$ws.Range("N122").Value = -11681.8
$ws.Range("H132").Value = 863661.2
$ws.Range("I132").Value = 927367.75
$ws.Range("J132").Value = 3622.5
$ws.Range("K132").Value = 2782103.25
$ws.Range("L132").Value = 10867.5
$ws.Range("M132").Value = -2779573.25
$ws.Range("N132").Value = -15927.5
$ws.Range("H135").Value = 89999
$ws.Range("J135").Value = 89999
$ws.Range("L135").Value = 89999
$ws.Range("N135").Value = -100139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 18042.75
$ws.Range("J100").Value = 18042.75
$ws.Range("L100").Value = 18042.75
$ws.Range("N100").Value = -20206.75
$ws.Range("H134").Value = 8374.826999999999
$ws.Range("I134").Value = 5102.5
$ws.Range("J134").Value = 100000
$ws.Range("K134").Value = 15307.5
$ws.Range("L134").Value = 300000
$ws.Range("M134").Value = -12772.5
$ws.Range("N134").Value = -305070
$ws.Range("H135").Value = 89249.5
$ws.Range("J135").Value = 89249.5
$ws.Range("L135").Value = 89249.5
$ws.Range("N135").Value = -99389.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 59.285713
$ws.Range("I7").Value = 52.5
$ws.Range("K7").Value = 52.5
$ws.Range("M7").Value = 60.5
$ws.Range("H11").Value = 1299
$ws.Range("I11").Value = 1248.75
$ws.Range("K11").Value = 1248.75
$ws.Range("M11").Value = -1108.75
$ws.Range("H16").Value = 1277.5238
$ws.Range("I16").Value = 1311.2778
$ws.Range("J16").Value = 1075
$ws.Range("K16").Value = 1311.2778
$ws.Range("L16").Value = 1075
$ws.Range("M16").Value = -1024.2778
$ws.Range("N16").Value = -1649
$ws.Range("H19").Value = 409.91666
$ws.Range("I19").Value = 438.0909
$ws.Range("K19").Value = 438.0909
$ws.Range("M19").Value = -268.0909
$ws.Range("H24").Value = 409.91666
$ws.Range("I24").Value = 438.0909
$ws.Range("K24").Value = 438.0909
$ws.Range("M24").Value = -268.0909
$ws.Range("H31").Value = 3670.838
$ws.Range("I31").Value = 4226.3335
$ws.Range("J31").Value = 2941.75
$ws.Range("K31").Value = 4226.3335
$ws.Range("L31").Value = 2941.75
$ws.Range("M31").Value = -3931.3335
$ws.Range("N31").Value = -3531.75
$ws.Range("H34").Value = 3670.838
$ws.Range("I34").Value = 4226.3335
$ws.Range("J34").Value = 2941.75
$ws.Range("K34").Value = 4226.3335
$ws.Range("L34").Value = 2941.75
$ws.Range("M34").Value = -4024.3335
$ws.Range("N34").Value = -3345.75
$ws.Range("H58").Value = 16198.5
$ws.Range("I58").Value = 6865.3335
$ws.Range("J58").Value = 30198.25
$ws.Range("K58").Value = 6865.3335
$ws.Range("L58").Value = 30198.25
$ws.Range("M58").Value = -6662.3335
$ws.Range("N58").Value = -30604.25
$ws.Range("H107").Value = 358.78946
$ws.Range("I107").Value = 307.3125
$ws.Range("K107").Value = 307.3125
$ws.Range("M107").Value = 1612.6875
$ws.Range("H113").Value = 1277.5238
$ws.Range("I113").Value = 1311.2778
$ws.Range("J113").Value = 1075
$ws.Range("K113").Value = 1311.2778
$ws.Range("L113").Value = 1075
$ws.Range("M113").Value = 858.7221999999999
$ws.Range("N113").Value = -5415
$ws.Range("H136").Value = 16198.5
$ws.Range("I136").Value = 6865.3335
$ws.Range("J136").Value = 30198.25
$ws.Range("K136").Value = 20596.0005
$ws.Range("L136").Value = 90594.75
$ws.Range("M136").Value = -18046.0005
$ws.Range("N136").Value = -95694.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1892325.8
$ws.Range("J4").Value = 3110.4211
$ws.Range("L4").Value = 9331.263300000001
$ws.Range("N4").Value = -9555.263300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 505000
$ws.Range("J7").Value = 9000
$ws.Range("L7").Value = 9000
$ws.Range("N7").Value = -9224
$ws.Range("H8").Value = 505000
$ws.Range("J8").Value = 9000
$ws.Range("L8").Value = 9000
$ws.Range("N8").Value = -9278
$ws.Range("H12").Value = 6000.5
$ws.Range("J12").Value = 9004
$ws.Range("L12").Value = 9004
$ws.Range("N12").Value = -9284
$ws.Range("H55").Value = 5775.4287
$ws.Range("I55").Value = 810
$ws.Range("J55").Value = 9499.5
$ws.Range("K55").Value = 810
$ws.Range("L55").Value = 9499.5
$ws.Range("M55").Value = -483
$ws.Range("N55").Value = -10153.5
$ws.Range("H92").Value = 18540.666
$ws.Range("J92").Value = 18540.666
$ws.Range("L92").Value = 18540.666
$ws.Range("N92").Value = -22284.666
$ws.Range("H107").Value = 4471.6665
$ws.Range("I107").Value = 6416.706
$ws.Range("K107").Value = 6416.706
$ws.Range("M107").Value = -4496.706
$ws.Range("H122").Value = 119656.555
$ws.Range("I122").Value = 252728
$ws.Range("K122").Value = 758184
$ws.Range("M122").Value = -755734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 12062.412
$ws.Range("I61").Value = 10378.9375
$ws.Range("K61").Value = 10378.9375
$ws.Range("M61").Value = -10176.9375
$ws.Range("H113").Value = 12062.412
$ws.Range("I113").Value = 10378.9375
$ws.Range("K113").Value = 10378.9375
$ws.Range("M113").Value = -8208.9375
$ws.Range("H122").Value = 3424.8572
$ws.Range("I122").Value = 3242.25
$ws.Range("J122").Value = 3668.3333
$ws.Range("K122").Value = 9726.75
$ws.Range("L122").Value = 11004.9999
$ws.Range("M122").Value = -7276.75
$ws.Range("N122").Value = -15904.9999
$ws.Range("H132").Value = 2283763.5
$ws.Range("I132").Value = 3586616.2
$ws.Range("J132").Value = 3771.25
$ws.Range("K132").Value = 10759848.6
$ws.Range("L132").Value = 11313.75
$ws.Range("M132").Value = -10757318.6
$ws.Range("N132").Value = -16373.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 135.90909
$ws.Range("I2").Value = 135.90909
$ws.Range("K2").Value = 135.90909
$ws.Range("M2").Value = -23.90908999999999
$ws.Range("H4").Value = 444899.12
$ws.Range("I4").Value = 571784.7
$ws.Range("K4").Value = 571784.7
$ws.Range("M4").Value = -571671.7
$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10336
